# ---------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" worksheet right after "总计" (i.e. right
#    before the existing "2022-Q2" sheet) and populate it with the
#    fund-holding detail rows for the new quarter.
# 2) Insert a new row at the top of the "总计" (totals) sheet's data
#    table for "2022-Q3", pushing the existing 2022-Q2 / 2022-Q1 /
#    2021-Q4 rows down by one.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet, positioned right before "2022-Q2"
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# Copy header row + index-column formatting/values from the 2022-Q2
# sheet so the new sheet matches the look of the other quarterly
# sheets (bold + bordered header row, bold + bordered index column).
$q2.Range("A1:H1").Copy($newSheet.Range("A1:H1"))
$q2.Range("A2:A6").Copy($newSheet.Range("A2:A6"))

# Fund code (B) and fund name (C) -- plain text columns.
$newSheet.Range("B2").Value = "001150"
$newSheet.Range("C2").Value = "融通互联网传媒灵活配置混合"
$newSheet.Range("B3").Value = "004818"
$newSheet.Range("C3").Value = "国寿安保目标策略灵活配置混合A"
$newSheet.Range("B4").Value = "004819"
$newSheet.Range("C4").Value = "国寿安保目标策略灵活配置混合C"
$newSheet.Range("B5").Value = "001834"
$newSheet.Range("C5").Value = "长盛战略新兴产业灵活配置混合C"
$newSheet.Range("B6").Value = "080008"
$newSheet.Range("C6").Value = "长盛战略新兴产业灵活配置混合A"

# Fund scale (D), stock position (E), position ratio (F) and held
# market value (G) are stored as text (numeric-looking strings, with
# significant trailing zeros) -- use a leading apostrophe so they are
# entered as text rather than being auto-converted to numbers.
$newSheet.Range("D2").Formula = "'8.61"
$newSheet.Range("E2").Formula = "'88.88"
$newSheet.Range("F2").Formula = "'2.68"
$newSheet.Range("G2").Formula = "'0.2307"

$newSheet.Range("D3").Formula = "'2.70"
$newSheet.Range("E3").Formula = "'45.00"
$newSheet.Range("F3").Formula = "'1.89"
$newSheet.Range("G3").Formula = "'0.0510"

$newSheet.Range("D4").Formula = "'1.73"
$newSheet.Range("E4").Formula = "'45.00"
$newSheet.Range("F4").Formula = "'1.89"
$newSheet.Range("G4").Formula = "'0.0327"

$newSheet.Range("D5").Formula = "'0.79"
$newSheet.Range("E5").Formula = "'50.54"
$newSheet.Range("F5").Formula = "'1.34"
$newSheet.Range("G5").Formula = "'0.0106"

$newSheet.Range("D6").Formula = "'0.12"
$newSheet.Range("E6").Formula = "'50.54"
$newSheet.Range("F6").Formula = "'1.34"
$newSheet.Range("G6").Formula = "'0.0016"

# Position rank (H) -- a real number column.
$newSheet.Range("H2").Value = 10
$newSheet.Range("H3").Value = 10
$newSheet.Range("H4").Value = 10
$newSheet.Range("H5").Value = 9
$newSheet.Range("H6").Value = 9

# ---------------------------------------------------------------
# 2) Shift the "总计" summary table down one row and fill in the new
#    2022-Q3 totals in the freed-up row 2.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push existing rows 2-4 down to rows 3-5, carrying formatting along
# (column A keeps its bold/bordered index style).
$total.Range("A4").Copy($total.Range("A5"))
$total.Range("B2:D4").Copy($total.Range("B3:D5"))

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.02

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 12
$total.Range("D4").Value = 1.31

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 26
$total.Range("D3").Value = 2.37

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.33

Write-Output "applied 2022-Q3 edits"
